# "Hoja con parámetros estadísticos del pronóstico"
# Adds an ETS.STAT statistics block (alfa, beta, gamma, MASE, SMAPE, MAE, ECM,
# step size) below the existing forecast block on the first sheet, and moves
# the active tab/selection from the "Estacionalidad" sheet back to the
# "Historico población" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Historico población"

# --- New statistics block: rows 35-42, columns A:C -----------------------

# Row 35: section header "ETS.STAT" (bold) + "Alfa" label + value
$ws1.Range("A35").Value = "ETS.STAT"
$ws1.Range("A35").Font.Bold = $true

$ws1.Range("B35").Value = "Alfa"
$ws1.Range("C35").NumberFormat = "0.000"
$ws1.Range("C35").Formula = "=_xlfn.FORECAST.ETS.STAT(B8:B21,A8:A21,1)"

# Row 36: Beta
$ws1.Range("B36").Value = "Beta"
$ws1.Range("C36").NumberFormat = "0.000"
$ws1.Range("C36").Formula = "=_xlfn.FORECAST.ETS.STAT(B8:B21,A8:A21,2)"

# Row 37: Gamma (very small value -> scientific notation)
$ws1.Range("B37").Value = "Gamma"
$ws1.Range("C37").NumberFormat = "0.000E+00"
$ws1.Range("C37").Formula = "=_xlfn.FORECAST.ETS.STAT(B8:B21,A8:A21,3)"

# Row 38: MASE
$ws1.Range("B38").Value = "MASE"
$ws1.Range("C38").NumberFormat = "0.000"
$ws1.Range("C38").Formula = "=_xlfn.FORECAST.ETS.STAT(B8:B21,A8:A21,4)"

# Row 39: SMAPE
$ws1.Range("B39").Value = "SMAPE"
$ws1.Range("C39").NumberFormat = "0.000"
$ws1.Range("C39").Formula = "=_xlfn.FORECAST.ETS.STAT(B8:B21,A8:A21,5)"

# Row 40: MAE
$ws1.Range("B40").Value = "MAE"
$ws1.Range("C40").NumberFormat = "0.000"
$ws1.Range("C40").Formula = "=_xlfn.FORECAST.ETS.STAT(B8:B21,A8:A21,6)"

# Row 41: ECM
$ws1.Range("B41").Value = "ECM"
$ws1.Range("C41").NumberFormat = "0.000"
$ws1.Range("C41").Formula = "=_xlfn.FORECAST.ETS.STAT(B8:B21,A8:A21,7)"

# Row 42: Tamaño de paso
$ws1.Range("B42").Value = "Tamaño de paso"
$ws1.Range("C42").NumberFormat = "0.000"
$ws1.Range("C42").Formula = "=_xlfn.FORECAST.ETS.STAT(B8:B21,A8:A21,8)"

# --- Page setup on sheet 1 (portrait, A4-ish "paperSize 9" = A4) ----------
$ps1 = $ws1.PageSetup
$ps1.PaperSize = 9
$ps1.Orientation = 1

# --- View / selection changes ---------------------------------------------
# Previously "Estacionalidad" (sheet 2) was the active tab with cell A20
# selected; now "Historico población" (sheet 1) becomes the active tab,
# scrolled down near the new statistics block, with C35 selected.
$ws1.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("C35").Select() | Out-Null
